$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new field names -----------------------------
$ws.Range("A1").Value = "firstName"
$ws.Range("B1").Value = "lastName"
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "phoneNumber"
$ws.Range("E1").Value = "zipCode"
$ws.Range("F1").ClearContents()

# New font for the header cells (A1:E1): Courier New, ~9.8pt, purple
$hdr = $ws.Range("A1:E1")
$hdr.Font.Name = "Courier New"
$hdr.Font.Size = 9.8
$hdr.Font.Color = 12287431

# --- Row 2: first contact entry ---------------------------------------
$ws.Range("A2").Value = "sample"
$ws.Range("B2").Value = "sample"
$ws.Range("C2").Value = "sample"
$ws.Range("D2").Value = 99293992
$ws.Range("E2").Value = 400101
$ws.Range("F2").ClearContents()

# --- Row 3: second contact entry (new) --------------------------------
$ws.Range("B3").Value = "sample1"
$ws.Range("C3").Value = "sample1"
$ws.Range("D3").Value = 289329399
$ws.Range("E3").Value = 400101

$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:wfsd2@gmail.com", "", "", "wfsd2@gmail.com")
$ws.Range("A3").Style = $ws.Range("A2").Style
$ws.Range("A3").Value = "sample1"

# --- Selection ---------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
